{"js": "const replacements = [\n  [\"2025-04-27 Sunday\", \"2025-04-28 Monday\"],\n  [\"60\u00d766=3960\", \"67\u00d738=2546\"],\n  [\"25\u00d730=750\", \"51\u00d757=2907\"],\n  [\"39\u00d753=2067\", \"15\u00d741=615\"],\n  [\"90\u00d721=1890\", \"41\u00d725=1025\"],\n  [\"15\u00d763=945\", \"36\u00d745=1620\"],\n  [\"33\u00d782=2706\", \"19\u00d784=1596\"],\n  [\"81\u00d764=5184\", \"99\u00d723=2277\"],\n  [\"52\u00d795=4940\", \"67\u00d731=2077\"],\n  [\"68\u00d793=6324\", \"47\u00d758=2726\"],\n  [\"46\u00d765=2990\", \"95\u00d795=9025\"],\n  [\"24\u00d772=1728\", \"56\u00d715=840\"],\n  [\"94\u00d714=1316\", \"99\u00d730=2970\"],\n  [\"21\u00d761=1281\", \"74\u00d791=6734\"],\n  [\"21\u00d759=1239\", \"52\u00d755=2860\"],\n  [\"24\u00d764=1536\", \"51\u00d751=2601\"],\n  [\"66\u00d790=5940\", \"58\u00d722=1276\"],\n  [\"39\u00d740=1560\", \"47\u00d785=3995\"],\n  [\"25\u00d714=350\", \"81\u00d751=4131\"],\n  [\"54\u00d762=3348\", \"51\u00d755=2805\"],\n  [\"97\u00d782=7954\", \"51\u00d751=2601\"],\n  [\"32\u00d775=2400\", \"79\u00d760=4740\"],\n  [\"77\u00d782=6314\", \"87\u00d765=5655\"],\n  [\"59\u00d775=4425\", \"15\u00d788=1320\"],\n  [\"72\u00d733=2376\", \"72\u00d785=6120\"],\n  [\"62\u00d727=1674\", \"61\u00d724=1464\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-04-27 Sunday\", \"2025-04-28 Monday\"),\n    @(\"60\u00d766=3960\", \"67\u00d738=2546\"),\n    @(\"25\u00d730=750\", \"51\u00d757=2907\"),\n    @(\"39\u00d753=2067\", \"15\u00d741=615\"),\n    @(\"90\u00d721=1890\", \"41\u00d725=1025\"),\n    @(\"15\u00d763=945\", \"36\u00d745=1620\"),\n    @(\"33\u00d782=2706\", \"19\u00d784=1596\"),\n    @(\"81\u00d764=5184\", \"99\u00d723=2277\"),\n    @(\"52\u00d795=4940\", \"67\u00d731=2077\"),\n    @(\"68\u00d793=6324\", \"47\u00d758=2726\"),\n    @(\"46\u00d765=2990\", \"95\u00d795=9025\"),\n    @(\"24\u00d772=1728\", \"56\u00d715=840\"),\n    @(\"94\u00d714=1316\", \"99\u00d730=2970\"),\n    @(\"21\u00d761=1281\", \"74\u00d791=6734\"),\n    @(\"21\u00d759=1239\", \"52\u00d755=2860\"),\n    @(\"24\u00d764=1536\", \"51\u00d751=2601\"),\n    @(\"66\u00d790=5940\", \"58\u00d722=1276\"),\n    @(\"39\u00d740=1560\", \"47\u00d785=3995\"),\n    @(\"25\u00d714=350\", \"81\u00d751=4131\"),\n    @(\"54\u00d762=3348\", \"51\u00d755=2805\"),\n    @(\"97\u00d782=7954\", \"51\u00d751=2601\"),\n    @(\"32\u00d775=2400\", \"79\u00d760=4740\"),\n    @(\"77\u00d782=6314\", \"87\u00d765=5655\"),\n    @(\"59\u00d775=4425\", \"15\u00d788=1320\"),\n    @(\"72\u00d733=2376\", \"72\u00d785=6120\"),\n    @(\"62\u00d727=1674\", \"61\u00d724=1464\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $found = $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Replacement not found for: $oldText\"\n    }\n}\n\nWrite-Output \"done\""}
